$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The very first paragraph in the document reads "*TODOS" (a yellow
# highlighted "*" run followed by a yellow highlighted "TODOS" run).
# A new yellow-highlighted run containing "SEG" must be inserted right
# before the "*" run, turning the paragraph into "SEG*TODOS".
$insertPoint = $d.Range(0, 0)
$insertPoint.InsertBefore("SEG")

# Re-select just the freshly inserted "SEG" text (characters 0-2) and make
# sure it carries the same yellow highlight as the run that follows it.
$segRange = $d.Range(0, 3)
$segRange.HighlightColorIndex = 7  # wdYellow

# --- Change 2 -----------------------------------------------------------
# Shorten the module heading from "Módulo: Catálogo y Recomendaciones" to
# "Módulo: Catálogo " (trailing space preserved).
$d.Content.Find.Execute("Módulo: Catálogo y Recomendaciones", $true, $false,
                         $false, $false, $false, $true, 1, $false,
                         "Módulo: Catálogo ", 2)
